# Applies the "stuff at the bottom of the sheets" commit:
#  1. Fills in the previously-blank pair_kind column (J2:J5) with "generic"
#     for the four practice pairs.
#  2. Appends a new "stim details" block starting at row 27 with its own
#     header row and a handful of video/audio stimulus-count rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- pair_kind for the practice rows (J2:J5) ---
$ws.Range("J2:J5").Value = "generic"

# --- new "stim details" section ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$stimRows = @(
    @{ Row = 29; Count = 6; Type = "video" },
    @{ Row = 30; Count = 6; Type = "video" },
    @{ Row = 31; Count = 7; Type = "video" },
    @{ Row = 32; Count = 7; Type = "video" },
    @{ Row = 33; Count = 6; Type = "audio" },
    @{ Row = 34; Count = 6; Type = "audio" },
    @{ Row = 35; Count = 7; Type = "audio" },
    @{ Row = 36; Count = 7; Type = "audio" }
)

foreach ($r in $stimRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Count
    $ws.Cells.Item($r.Row, 2).Value = $r.Type
}
